$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("MISC", ".github/workflows/testPR.yml"),
    @("MISC", "Impacted_Files_List.xlsx"),
    @("MISC", "Impacted_Modules_List.xlsx"),
    @("MISC", "test.txt"),
    @("RN_LIBRARIES", "package.json"),
    @("RN_LIBRARIES", "yarn.lock")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
